# Weekly update: insert a new week's worth of "Plátano" price rows at the
# top of the data block (row 611), pushing the existing rows (previously
# 611-665) down by 4, to 615-669.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows before the current row 611. Excel COM copies the
# formatting of the row above the insertion point, so column D keeps its
# date-formatted style (index 2), matching the rest of the sheet.
$ws.Rows("611:614").Insert()

# Row 611: Barraganete / Primera
$ws.Cells.Item(611, 1).Value2  = 9
$ws.Cells.Item(611, 2).Value2  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(611, 3).Value2  = "Metropolitana"
$ws.Cells.Item(611, 4).Value2  = 44461
$ws.Cells.Item(611, 5).Value2  = 13
$ws.Cells.Item(611, 6).Value2  = "Fruta"
$ws.Cells.Item(611, 7).Value2  = 100108
$ws.Cells.Item(611, 8).Value2  = "Tropicales y subtropicales"
$ws.Cells.Item(611, 9).Value2  = 100108006
$ws.Cells.Item(611, 10).Value2 = "Plátano"
$ws.Cells.Item(611, 11).Value2 = "Barraganete"
$ws.Cells.Item(611, 12).Value2 = "Primera"
$ws.Cells.Item(611, 13).Value2 = 220
$ws.Cells.Item(611, 14).Value2 = 18000
$ws.Cells.Item(611, 15).Value2 = 19000
$ws.Cells.Item(611, 16).Value2 = 18455
$ws.Cells.Item(611, 17).Value2 = "$/caja 20 kilos"
$ws.Cells.Item(611, 18).Value2 = "Ecuador"
$ws.Cells.Item(611, 19).Value2 = 923
$ws.Cells.Item(611, 20).Value2 = 20

# Row 612: Sin especificar / Pintón
$ws.Cells.Item(612, 1).Value2  = 9
$ws.Cells.Item(612, 2).Value2  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(612, 3).Value2  = "Metropolitana"
$ws.Cells.Item(612, 4).Value2  = 44461
$ws.Cells.Item(612, 5).Value2  = 13
$ws.Cells.Item(612, 6).Value2  = "Fruta"
$ws.Cells.Item(612, 7).Value2  = 100108
$ws.Cells.Item(612, 8).Value2  = "Tropicales y subtropicales"
$ws.Cells.Item(612, 9).Value2  = 100108006
$ws.Cells.Item(612, 10).Value2 = "Plátano"
$ws.Cells.Item(612, 11).Value2 = "Sin especificar"
$ws.Cells.Item(612, 12).Value2 = "Pintón"
$ws.Cells.Item(612, 13).Value2 = 320
$ws.Cells.Item(612, 14).Value2 = 10000
$ws.Cells.Item(612, 15).Value2 = 11000
$ws.Cells.Item(612, 16).Value2 = 10562
$ws.Cells.Item(612, 17).Value2 = "$/caja 20 kilos"
$ws.Cells.Item(612, 18).Value2 = "Ecuador"
$ws.Cells.Item(612, 19).Value2 = 528
$ws.Cells.Item(612, 20).Value2 = 20

# Row 613: Sin especificar / Primera Maduro
$ws.Cells.Item(613, 1).Value2  = 9
$ws.Cells.Item(613, 2).Value2  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(613, 3).Value2  = "Metropolitana"
$ws.Cells.Item(613, 4).Value2  = 44461
$ws.Cells.Item(613, 5).Value2  = 13
$ws.Cells.Item(613, 6).Value2  = "Fruta"
$ws.Cells.Item(613, 7).Value2  = 100108
$ws.Cells.Item(613, 8).Value2  = "Tropicales y subtropicales"
$ws.Cells.Item(613, 9).Value2  = 100108006
$ws.Cells.Item(613, 10).Value2 = "Plátano"
$ws.Cells.Item(613, 11).Value2 = "Sin especificar"
$ws.Cells.Item(613, 12).Value2 = "Primera Maduro"
$ws.Cells.Item(613, 13).Value2 = 1180
$ws.Cells.Item(613, 14).Value2 = 12000
$ws.Cells.Item(613, 15).Value2 = 13000
$ws.Cells.Item(613, 16).Value2 = 12492
$ws.Cells.Item(613, 17).Value2 = "$/caja 20 kilos"
$ws.Cells.Item(613, 18).Value2 = "Ecuador"
$ws.Cells.Item(613, 19).Value2 = 625
$ws.Cells.Item(613, 20).Value2 = 20

# Row 614: Sin especificar / Primera Pintón
$ws.Cells.Item(614, 1).Value2  = 9
$ws.Cells.Item(614, 2).Value2  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(614, 3).Value2  = "Metropolitana"
$ws.Cells.Item(614, 4).Value2  = 44461
$ws.Cells.Item(614, 5).Value2  = 13
$ws.Cells.Item(614, 6).Value2  = "Fruta"
$ws.Cells.Item(614, 7).Value2  = 100108
$ws.Cells.Item(614, 8).Value2  = "Tropicales y subtropicales"
$ws.Cells.Item(614, 9).Value2  = 100108006
$ws.Cells.Item(614, 10).Value2 = "Plátano"
$ws.Cells.Item(614, 11).Value2 = "Sin especificar"
$ws.Cells.Item(614, 12).Value2 = "Primera Pintón"
$ws.Cells.Item(614, 13).Value2 = 1040
$ws.Cells.Item(614, 14).Value2 = 13000
$ws.Cells.Item(614, 15).Value2 = 14000
$ws.Cells.Item(614, 16).Value2 = 13538
$ws.Cells.Item(614, 17).Value2 = "$/caja 20 kilos"
$ws.Cells.Item(614, 18).Value2 = "Ecuador"
$ws.Cells.Item(614, 19).Value2 = 677
$ws.Cells.Item(614, 20).Value2 = 20
